$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# about_main (B3): "6 concepts" -> "7 concepts"
$txtB3 = @'
![main_banner not_rounded](data/img/main_banner{dark_mode}.png?v=1)

# datannur, le catalogue de données portable

Permet de **centraliser**, **rechercher** et **visualiser** les informations sur une collection de jeux de données

Pour améliorer l’organisation des données et faciliter leur **partage** et leur **documentation**

**Simple** et **flexible**, s’intègre rapidement dans tous types d’environnement


- **Facile** :
Aucune installation ou configuration nécessaire, aucun coût ou prérequis technique

- **Portable** :
Fonctionne partout (local, cloud, disque partagé), un simple dossier que l’on peut copier, déplacer, envoyer et ouvrir avec n’importe quel navigateur web

- **Complet** :
Flexible, complet et structuré autour de 7 concepts avec un niveau de détail important : Institution, Dossier, Mot clé, Dataset, Variable et Modalité

- **Indépendant** :
Le catalogue n’est qu’une interface pour visualiser les métadonnées, le processus de leur création et mise à jour est indépendant et sous votre contrôle

- **Sécurisé** :
De pars la séparation stricte entre les deux systèmes, l’application est isolée dans le navigateur, ne peut rien modifier sur la machine et ne pose ainsi aucun risque

La version ici présente est un **prototype** en cours de développement et d'expérimentation. Les données utilisées sont fictives et uniquement à usage de test et de développement. Question ou suggestion : [contact@datannur.com](mailto:contact@datannur.com).
'@
$ws.Range("B3").Value = $txtB3

# about_page_1 (B4): "7 entités principales" -> "7 concepts principaux"
$txtB4 = @'
### Fonctionnement
datannur contient 7 concepts principaux. On peut les diviser en deux catégories, partie intérieur et partie extérieur aux datasets. Le **dataset** représente une table de base de données ou un fichier de données (excel, csv, ...) sous forme de tableau (lignes et colonnes).

mermaid( 
  $dataset -.-> intérieur
  $dataset -.-> extérieur
);

'@
$ws.Range("B4").Value = $txtB4

# about_page_4 (B7): "certaines entités" -> "certains concepts"
$txtB7 = @'
Pour finir, certains concepts possèdent des docs (documents de type markdown ou pdf). Un **doc** peut être lié à une multitude de datasets, de dossiers et d'institutions, et inversement.

mermaid(
  $doc <--> $institution
  $doc <--> $folder
  $doc <--> $dataset
);
'@
$ws.Range("B7").Value = $txtB7

# about_page_5 (B8): "7 entités" -> "7 concepts"
$txtB8 = @'
#### Vision d'ensemble

Voici les liens entre les 7 concepts, à l'intérieur et l'extérieur des datasets.

mermaid(
  $folder $recursive
  $institution $recursive
  $tag $recursive
  $institution -- manager - owner --> $dataset
  $institution -- manager - owner --> $folder
  $folder --> $dataset
  $folder --> $modality
  $tag <--> $institution
  $tag <--> $folder
  $tag <--> $dataset
  $tag <--> $variable
  $doc <--> $institution
  $doc <--> $folder
  $doc <--> $dataset
  $dataset --> $variable
  $variable <--> $modality
  $modality --> $value
  );

'@
$ws.Range("B8").Value = $txtB8

# about_doc (B18): "Certaines entités" -> "Certains concepts"
$txtB18 = @'
Certains concepts possèdent des docs (documents de type markdown ou pdf). Un **doc** peut être lié à une multitude de datasets, de dossiers et d'institutions, et inversement.
'@
$ws.Range("B18").Value = $txtB18

# Restore the saved selection/scroll state recorded in the workbook.
$null = $ws.Range("B22").Select()
